# Update "想去人数" (interested-count) figures to the latest scraped values.
# Source data refresh corresponds to commit 456a3b4 (gh-pages generated output).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value  = 123
$wsExpo.Range("F8").Value  = 4666
$wsExpo.Range("F9").Value  = 101
$wsExpo.Range("F10").Value = 5074
$wsExpo.Range("F11").Value = 580
$wsExpo.Range("F12").Value = 1271

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value   = 123
$wsAll.Range("F9").Value   = 4666
$wsAll.Range("F10").Value  = 101
$wsAll.Range("F11").Value  = 5074
$wsAll.Range("F12").Value  = 580
$wsAll.Range("F13").Value  = 1271
